# Refresh cryptos list: update Price (D) and Volume(1h) (E) columns
# to the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.610.62'
$ws.Range('E2').Value = '  -0.33%  '
$ws.Range('D3').Value = '1.596.75'
$ws.Range('E3').Value = '  -0.18%  '
$ws.Range('D5').Value = '''210.56'
$ws.Range('E5').Value = '  -0.44%  '
$ws.Range('D6').Value = '''0.509'
$ws.Range('E6').Value = '  -0.70%  '
$ws.Range('E7').Value = '  +0.38%  '
$ws.Range('D8').Value = '''0.0614'
$ws.Range('E8').Value = '  -0.78%  '
$ws.Range('D9').Value = '''0.245'
$ws.Range('E9').Value = '  -0.82%  '
$ws.Range('D10').Value = '''19.54'
$ws.Range('E10').Value = '  -0.08%  '
$ws.Range('D11').Value = '''0.0845'
$ws.Range('E11').Value = '  +0.48%  '
$ws.Range('D12').Value = '1.821.40'
$ws.Range('E12').Value = '  -0.16%  '
$ws.Range('D13').Value = '1.604.22'
$ws.Range('E13').Value = '  +0.45%  '
$ws.Range('E14').Value = '  -0.12%  '
$ws.Range('E15').Value = '  -0.38%  '
$ws.Range('D16').Value = '''64.63'
$ws.Range('E16').Value = '  -1.06%  '
$ws.Range('D17').Value = '26.606.44'
$ws.Range('E17').Value = '  -0.27%  '
$ws.Range('D18').Value = '0.0₃0737'
$ws.Range('E18').Value = '  -2.64%  '
$ws.Range('E19').Value = '  +0.29%  '
$ws.Range('D20').Value = '''208.36'
$ws.Range('E20').Value = '  -0.73%  '
$ws.Range('E21').Value = '  -1.17%  '
$ws.Range('E22').Value = '  -0.28%  '
$ws.Range('D23').Value = '''2.23'
$ws.Range('E23').Value = '  -3.38%  '
$ws.Range('D24').Value = '''8.91'
$ws.Range('E24').Value = '  -0.17%  '
$ws.Range('D25').Value = '''143.76'
$ws.Range('E25').Value = '  +0.62%  '
$ws.Range('E26').Value = '  +0.37%  '
$ws.Range('E27').Value = '  +0.09%  '
$ws.Range('E28').Value = '  -1.15%  '
$ws.Range('E29').Value = '  -0.61%  '
$ws.Range('D30').Value = '''0.0506'
$ws.Range('E30').Value = '  -2.19%  '
$ws.Range('E31').Value = '  -0.21%  '
$ws.Range('E32').Value = '  -0.61%  '
$ws.Range('E33').Value = '  -0.58%  '
$ws.Range('E34').Value = '  +19.09%  '
$ws.Range('D35').Value = '1.275.40'
$ws.Range('E35').Value = '  -1.19%  '
$ws.Range('E36').Value = '  +0.92%  '
$ws.Range('E37').Value = '  -0.85%  '
$ws.Range('D38').Value = '''0.596'
$ws.Range('E38').Value = '  -3.60%  '
$ws.Range('E39').Value = '  -2.18%  '
$ws.Range('E40').Value = '  -0.67%  '
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('E42').Value = '  -1.29%  '
$ws.Range('D44').Value = '''62.48'
$ws.Range('E44').Value = '  -1.06%  '
$ws.Range('D45').Value = '1.733.19'
$ws.Range('E45').Value = '  -0.12%  '
$ws.Range('D46').Value = '''89.74'
$ws.Range('E46').Value = '  -1.46%  '
$ws.Range('E47').Value = '  -0.65%  '
$ws.Range('E48').Value = '  -2.10%  '
$ws.Range('E49').Value = '  +2.03%  '
$ws.Range('D50').Value = '''0.0513'
$ws.Range('E50').Value = '  +0.76%  '
$ws.Range('E51').Value = '  +0.30%  '
